# Asesores.xlsx — "Add files via upload"
#
# The underlying change is a single new advisor row inserted into the
# "ALHUANUCO" store group on sheet Hoja1: a brand-new row 97 carrying
# A97 = "ALHUANUCO" (same store code as the group it's appended to) and
# B97 = "PEÑA ESPINOZA JHERRISON MARIO" (a brand-new shared string).
# Every existing row from (old) row 97 downward shifts down by one, which
# Excel's native row-insert semantics handle automatically (formulas,
# dimension, etc. all re-anchor on their own).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert a new blank row above the first row of the ALCHINCHA group (old
# row 97), which lands it as the last row of the preceding ALHUANUCO
# group — exactly where the diff shows it.
$ws.Rows("97:97").Insert()

# Populate the new row: same store code as the rest of the ALHUANUCO
# block directly above it, plus the new advisor's name. (Read via
# .Value2 — .Value's getter doesn't resolve through this COM shim.)
$ws.Range("A97").Value = $ws.Range("A96").Value2
$ws.Range("B97").Value = "PEÑA ESPINOZA JHERRISON MARIO"

# Column A picks up an explicit best-fit width once it holds real data
# (mirrors the <col min="1" max="1" .../> that appears in the saved XML).
$ws.Columns("A:A").AutoFit()

# Cosmetic view-state nudges (selection / scroll anchor) to mirror the
# author's on-screen state after the edit.
$ws.Range("G97").Select()
